$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the figure file name and status for row 15 (item #14)
$ws.Range("C15").Value = "streaks.eps"
$ws.Range("D15").Value = "ok"

# Move the active selection to C16 (matches the saved sheetView selection)
$ws.Range("C16").Select()
